# Insert a new weekly price-report row at row 182 (shifting existing rows
# 182-213 down to 183-214) and populate it with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 182; rows 182:213 shift down to 183:214.
$ws.Rows(182).EntireRow.Insert()

$row = 182

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44522
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112032
$ws.Cells.Item($row, 7).Value = "Zapallo italiano"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 400
$ws.Cells.Item($row, 11).Value = 10000
$ws.Cells.Item($row, 12).Value = 11000
$ws.Cells.Item($row, 13).Value = 10500
$ws.Cells.Item($row, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item($row, 15).Value = "Provincia de Limarí"
$ws.Cells.Item($row, 16).Value = 150
$ws.Cells.Item($row, 17).Value = 70
$ws.Cells.Item($row, 18).Value = "Hortaliza"
